$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32: Automata for the People | Crab Oil
$ws.Range("H32").Value = 1366.6666
$ws.Range("J32").Value = 1539.4286
$ws.Range("L32").Value = 1539.4286
$ws.Range("N32").Value = -2191.4286

# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 33354048
$ws.Range("I40").Value = 25615.2
$ws.Range("J40").Value = 50018264
$ws.Range("K40").Value = 25615.2
$ws.Range("L40").Value = 50018264
$ws.Range("M40").Value = -25440.2
$ws.Range("N40").Value = -50018614

# Row 45: The House Always Wins | Blinding Potion
$ws.Range("H45").Value = 12000
$ws.Range("J45").Value = 12000
$ws.Range("L45").Value = 36000
$ws.Range("N45").Value = -36384

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 1380.138
$ws.Range("I98").Value = 1228.7727
$ws.Range("K98").Value = 1228.7727
$ws.Range("M98").Value = 269.2273

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 173616750
$ws.Range("I116").Value = 180556850
$ws.Range("J116").Value = 166676670
$ws.Range("K116").Value = 180556850
$ws.Range("L116").Value = 166676670
$ws.Range("M116").Value = -180553408
$ws.Range("N116").Value = -166683554

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 1380.138
$ws.Range("I122").Value = 1228.7727
$ws.Range("K122").Value = 3686.3181
$ws.Range("M122").Value = -1236.3181

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 5863.902
$ws.Range("J138").Value = 7265.675
$ws.Range("L138").Value = 21797.025
$ws.Range("N138").Value = -32077.025

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 4484.6665
$ws.Range("I141").Value = 4597.364
$ws.Range("J141").Value = 3988.8
$ws.Range("K141").Value = 13792.092
$ws.Range("L141").Value = 11966.4
$ws.Range("M141").Value = -8612.091999999999
$ws.Range("N141").Value = -22326.4

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth | Bronze Rivets
$ws.Range("H5").Value = 583
$ws.Range("I5").Value = 665.1667
$ws.Range("K5").Value = 665.1667
$ws.Range("M5").Value = -553.1667

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 4327.4287
$ws.Range("J45").Value = 4333.3335
$ws.Range("L45").Value = 4333.3335
$ws.Range("N45").Value = -5087.3335

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 27859.412
$ws.Range("I61").Value = 46876
$ws.Range("K61").Value = 46876
$ws.Range("M61").Value = -46664

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 5323117.5
$ws.Range("I74").Value = 10417777
$ws.Range("J74").Value = 6950.8696
$ws.Range("K74").Value = 10417777
$ws.Range("L74").Value = 6950.8696
$ws.Range("M74").Value = -10416903
$ws.Range("N74").Value = -8698.8696

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 5323117.5
$ws.Range("I77").Value = 10417777
$ws.Range("J77").Value = 6950.8696
$ws.Range("K77").Value = 52088885
$ws.Range("L77").Value = 34754.348
$ws.Range("M77").Value = -52084517
$ws.Range("N77").Value = -43490.348

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 2753
$ws.Range("I122").Value = 2753
$ws.Range("K122").Value = 8259
$ws.Range("M122").Value = -5809

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 19283.262
$ws.Range("I132").Value = 19316.908
$ws.Range("K132").Value = 57950.724
$ws.Range("M132").Value = -55420.724

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 27859.412
$ws.Range("I136").Value = 46876
$ws.Range("K136").Value = 140628
$ws.Range("M136").Value = -138078

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences | Bronze Rivets
$ws.Range("H4").Value = 583
$ws.Range("I4").Value = 665.1667
$ws.Range("K4").Value = 665.1667
$ws.Range("M4").Value = -550.1667

# Row 22: Riveting Run | Iron Rivets
$ws.Range("H22").Value = 426.13333
$ws.Range("I22").Value = 426.13333
$ws.Range("K22").Value = 426.13333
$ws.Range("M22").Value = -253.13333

# Row 52: File That under Whatever | Mythril File
$ws.Range("H52").Value = 55000
$ws.Range("J52").Value = 55000
$ws.Range("L52").Value = 55000
$ws.Range("N52").Value = -55526

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 78526920
$ws.Range("I99").Value = 92804320
$ws.Range("K99").Value = 92804320
$ws.Range("M99").Value = -92802822

# Row 121: Keeping Loyalty | Dwarven Mythril File
$ws.Range("H121").Value = 55000
$ws.Range("J121").Value = 55000
$ws.Range("L121").Value = 55000
$ws.Range("N121").Value = -58494

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 3745.3403
$ws.Range("I134").Value = 1076.5555
$ws.Range("K134").Value = 3229.6665
$ws.Range("M134").Value = -694.6664999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 21280958
$ws.Range("I31").Value = 62501936
$ws.Range("K31").Value = 62501936
$ws.Range("M31").Value = -62501641

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 21280958
$ws.Range("I34").Value = 62501936
$ws.Range("K34").Value = 62501936
$ws.Range("M34").Value = -62501734

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 3334724.8
$ws.Range("I58").Value = 5000937
$ws.Range("J58").Value = 2300
$ws.Range("K58").Value = 5000937
$ws.Range("L58").Value = 2300
$ws.Range("M58").Value = -5000734
$ws.Range("N58").Value = -2706

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 78445560
$ws.Range("I132").Value = 95240840
$ws.Range("J132").Value = 67604.336
$ws.Range("K132").Value = 285722520
$ws.Range("L132").Value = 202813.008
$ws.Range("M132").Value = -285719990
$ws.Range("N132").Value = -207873.008

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 2323.4443
$ws.Range("I134").Value = 1752.1666
$ws.Range("J134").Value = 3466
$ws.Range("K134").Value = 5256.4998
$ws.Range("L134").Value = 10398
$ws.Range("M134").Value = -2721.4998
$ws.Range("N134").Value = -15468

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 3334724.8
$ws.Range("I136").Value = 5000937
$ws.Range("J136").Value = 2300
$ws.Range("K136").Value = 15002811
$ws.Range("L136").Value = 6900
$ws.Range("M136").Value = -15000261
$ws.Range("N136").Value = -12000

$ws = $wb.Worksheets.Item("CUL")
# Row 44: No More Dumpster Diving | Knight's Bread
$ws.Range("H44").Value = 407
$ws.Range("J44").Value = 462.25
$ws.Range("L44").Value = 1386.75
$ws.Range("N44").Value = -2182.75

# Row 68: Such a Butter Face | Fermented Butter
$ws.Range("H68").Value = 180276.7
$ws.Range("I68").Value = 1999.75
$ws.Range("K68").Value = 5999.25
$ws.Range("M68").Value = -5188.25

# Row 71: No Margarine of Error (L) | Fermented Butter
$ws.Range("H71").Value = 180276.7
$ws.Range("I71").Value = 1999.75
$ws.Range("K71").Value = 17997.75
$ws.Range("M71").Value = -13941.75

# Row 113: Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 627.1429000000001
$ws.Range("I113").Value = 547.5
$ws.Range("J113").Value = 659
$ws.Range("K113").Value = 1642.5
$ws.Range("L113").Value = 1977
$ws.Range("M113").Value = 527.5
$ws.Range("N113").Value = -6317

# Row 137: Creative Chocolate | Gateau au Chocolat
$ws.Range("H137").Value = 43002910
$ws.Range("I137").Value = 44118720
$ws.Range("J137").Value = 36680000
$ws.Range("K137").Value = 132356160
$ws.Range("L137").Value = 110040000
$ws.Range("M137").Value = -132351060
$ws.Range("N137").Value = -110050200

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 2930.0278
$ws.Range("I132").Value = 2762.5
$ws.Range("K132").Value = 8287.5
$ws.Range("M132").Value = -5757.5

$ws = $wb.Worksheets.Item("LTW")
# Row 6: Sticking Their Necks Out | Leather Choker
$ws.Range("H6").Value = 31623.875
$ws.Range("J6").Value = 31623.875
$ws.Range("L6").Value = 31623.875
$ws.Range("N6").Value = -31847.875

# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 40409036
$ws.Range("I40").Value = 18522676
$ws.Range("J40").Value = 66672668
$ws.Range("K40").Value = 18522676
$ws.Range("L40").Value = 66672668
$ws.Range("M40").Value = -18522540
$ws.Range("N40").Value = -66672940

# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 2274037.2
$ws.Range("I68").Value = 3789028.8
$ws.Range("K68").Value = 3789028.8
$ws.Range("M68").Value = -3788279.8

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 2274037.2
$ws.Range("I71").Value = 3789028.8
$ws.Range("K71").Value = 18945144
$ws.Range("M71").Value = -18941400

# Row 82: Trainin' the Neck | Dragon Leather
$ws.Range("H82").Value = 2233678.5
$ws.Range("I82").Value = 3473133
$ws.Range("K82").Value = 3473133
$ws.Range("M82").Value = -3472772

# Row 85: Training Is Only Skintight (L) | Dragon Leather
$ws.Range("H85").Value = 2233678.5
$ws.Range("I85").Value = 3473133
$ws.Range("K85").Value = 3473133
$ws.Range("M85").Value = -3471885

# Row 118: Strike True | Zonureskin Fingerless Gloves of Aiming
$ws.Range("H118").Value = 11950
$ws.Range("J118").Value = 11950
$ws.Range("L118").Value = 11950
$ws.Range("N118").Value = -15264

# Row 127: Loyal Turncoat | Saigaskin Coat of Fending
$ws.Range("H127").Value = 3385666.8
$ws.Range("J127").Value = 78500
$ws.Range("L127").Value = 78500
$ws.Range("N127").Value = -88420

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 3980.1238
$ws.Range("I136").Value = 3979.5652
$ws.Range("J136").Value = 3990.4
$ws.Range("K136").Value = 11938.6956
$ws.Range("L136").Value = 11971.2
$ws.Range("M136").Value = -9388.695599999999
$ws.Range("N136").Value = -17071.2

$ws = $wb.Worksheets.Item("WVR")
# Row 45: Private Concerns | Linen Trousers
$ws.Range("H45").Value = 15392
$ws.Range("J45").Value = 11784
$ws.Range("L45").Value = 11784
$ws.Range("N45").Value = -12766

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 9807299
$ws.Range("I132").Value = 3642.5625
$ws.Range("K132").Value = 10927.6875
$ws.Range("M132").Value = -8397.6875

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 6244.0303
$ws.Range("I136").Value = 2885.1316
$ws.Range("J136").Value = 8336.459000000001
$ws.Range("K136").Value = 8655.3948
$ws.Range("L136").Value = 25009.377
$ws.Range("M136").Value = -6105.3948
$ws.Range("N136").Value = -30109.377
